$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "66.974.00"
$ws.Range("E2").Value = "  -0.20%  "

# Row 3
$ws.Range("D3").Value = "3.087.15"
$ws.Range("E3").Value = "  -1.33%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'578.92"
$ws.Range("E5").Value = "  +0.07%  "

# Row 6
$ws.Range("D6").Value = "'169.84"
$ws.Range("E6").Value = "  -2.51%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").Value = "3.082.76"
$ws.Range("E8").Value = "  -1.32%  "

# Row 9
$ws.Range("D9").Value = "'0.515"
$ws.Range("E9").Value = "  -1.64%  "

# Row 10
$ws.Range("E10").Value = "  -0.46%  "

# Row 11
$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = "  -3.13%  "

# Row 12
$ws.Range("D12").Value = "'0.474"
$ws.Range("E12").Value = "  -1.36%  "

# Row 13
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("E13").Value = "  -2.68%  "

# Row 14
$ws.Range("D14").Value = "'36.46"
$ws.Range("E14").Value = "  -2.16%  "

# Row 15
$ws.Range("E15").Value = "  -2.14%  "

# Row 16
$ws.Range("D16").Value = "3.598.06"
$ws.Range("E16").Value = "  -1.34%  "

# Row 17
$ws.Range("D17").Value = "66.919.17"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18
$ws.Range("D18").Value = "'7.05"
$ws.Range("E18").Value = "  -1.25%  "

# Row 19
$ws.Range("D19").Value = "3.083.57"
$ws.Range("E19").Value = "  -1.54%  "

# Row 20
$ws.Range("D20").Value = "'16.42"
$ws.Range("E20").Value = "  +1.62%  "

# Row 21
$ws.Range("D21").Value = "'484.36"
$ws.Range("E21").Value = "  -0.07%  "

# Row 22
$ws.Range("D22").Value = "'7.75"
$ws.Range("E22").Value = "  +0.42%  "

# Row 23
$ws.Range("D23").Value = "'0.692"
$ws.Range("E23").Value = "  -3.17%  "

# Row 24
$ws.Range("D24").Value = "'83.11"
$ws.Range("E24").Value = "  -1.06%  "

# Row 25
$ws.Range("D25").Value = "'12.91"
$ws.Range("E25").Value = "  -3.25%  "

# Row 26
$ws.Range("D26").Value = "'2.25"
$ws.Range("E26").Value = "  -3.11%  "

# Row 27
$ws.Range("D27").Value = "'10.36"
$ws.Range("E27").Value = "  +2.93%  "

# Row 28
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("E29").Value = "  -3.27%  "

# Row 30
$ws.Range("E30").Value = "  -3.55%  "

# Row 31
$ws.Range("D31").Value = "'2.64"
$ws.Range("E31").Value = "  -1.55%  "

# Row 32
$ws.Range("D32").Value = "'27.98"
$ws.Range("E32").Value = "  -2.85%  "

# Row 33
$ws.Range("E33").Value = "  -1.91%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0927"
$ws.Range("E34").Value = "  -6.50%  "

# Row 35
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.15%  "

# Row 36
$ws.Range("D36").Value = "'5.71"
$ws.Range("E36").Value = "  -2.97%  "

# Row 37
$ws.Range("D37").Value = "'0.959"
$ws.Range("E37").Value = "  -2.44%  "

# Row 38
$ws.Range("D38").Value = "'46.50"
$ws.Range("E38").Value = "  -2.43%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.124"
$ws.Range("E39").Value = "  +0.36%  "

# Row 40
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'1.99"
$ws.Range("E40").Value = "  -4.61%  "

# Row 41
$ws.Range("D41").Value = "'0.303"
$ws.Range("E41").Value = "  -2.93%  "

# Row 42
$ws.Range("D42").Value = "'8.36"
$ws.Range("E42").Value = "  -3.23%  "

# Row 43
$ws.Range("D43").Value = "2.777.80"
$ws.Range("E43").Value = "  -2.49%  "

# Row 44
$ws.Range("D44").Value = "'379.30"
$ws.Range("E44").Value = "  -0.89%  "

# Row 45
$ws.Range("D45").Value = "'2.55"
$ws.Range("E45").Value = "  -4.58%  "

# Row 46
$ws.Range("D46").Value = "'0.0347"
$ws.Range("E46").Value = "  -3.02%  "

# Row 47
$ws.Range("D47").Value = "'135.27"
$ws.Range("E47").Value = "  -0.57%  "

# Row 48
$ws.Range("E48").Value = "  +0.00%  "

# Row 49
$ws.Range("D49").Value = "'24.52"
$ws.Range("E49").Value = "  -1.53%  "

# Row 50
$ws.Range("D50").Value = "'2.16"
$ws.Range("E50").Value = "  -2.61%  "

# Row 51
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  -1.85%  "
